$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "IP" header cell (H1) onto the two
# new header cells so they pick up the same bold/centered/bordered style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New column data (I = I0, J = IF) for rows 2-14
$data = @(
    @(1, 5),
    @(5, 8),
    @(2, 4),
    @(3, 8),
    @(3, 5),
    @(1, 4),
    @(3, 4),
    @(2, 5),
    @(1, 4),
    @(5, 8),
    @(1, 4),
    @(1, 2),
    @(1, 1)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
